$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.747982217853126
$ws.Range("C2").Value = 1.272608002624571
$ws.Range("D2").Value = 0.3342653064198473

$ws.Range("B3").Value = 2.619778015857148
$ws.Range("C3").Value = 1.924682895761365
$ws.Range("D3").Value = 2.376979159805273

$ws.Range("B4").Value = 2.739185677413226
$ws.Range("C4").Value = -0.4644821689663494
$ws.Range("D4").Value = -0.6014901743245544

$ws.Range("B5").Value = 6.102051446073806
$ws.Range("C5").Value = 1.213248978582696
$ws.Range("D5").Value = 0.4328646475346514

$ws.Range("B6").Value = 4.850344730476794
$ws.Range("C6").Value = 0.05769090177856884
$ws.Range("D6").Value = 0.47118629539715

$ws.Range("B7").Value = 5.440657912325899
$ws.Range("C7").Value = 1.996251390219148
$ws.Range("D7").Value = 1.186194765167633
